# Auto-generated Excel COM-interop script
# Appends transaction rows 72-86 (2019-04-20 purchases) to the "Transacciones" sheet,
# mirroring the shared running-balance formulas used throughout the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

$dateSerial = 43575  # 2019-04-20, one day after the last existing entry (row 71)

# --- Row 72: Gasolina ---
$ws.Range("A72").Value = $dateSerial
$ws.Range("B72").Value = 163
$ws.Range("C72").Value = "Gasolina"
$ws.Range("D72").Value = "Gasolina"
$ws.Range("E72").Value = "Gasto"
$ws.Range("F72").Value = "Efectivo"
$ws.Range("G72").Value = "Gasolinería Mobil"
$ws.Range("K72").Value = 7900.24
$ws.Range("L72").Value = 2527.5700000000002
$ws.Range("M72").Formula = "=M71-B72"
$ws.Range("N72").Formula = "=SUM(K72:M72)"
$ws.Range("O72").Formula = "=N72-4000"

# --- Row 73: Atún Dolores en Agua ---
$ws.Range("A73").Value = $dateSerial
$ws.Range("B73").Value = 16.3
$ws.Range("C73").Value = "Atún Dolores en Agua"
$ws.Range("D73").Value = "Despensa"
$ws.Range("E73").Value = "Gasto"
$ws.Range("F73").Value = "Tarjeta Banamex"
$ws.Range("G73").Value = "Soriana"
$ws.Range("K73").Formula = "=K72-B73"
$ws.Range("L73").Value = 2527.5700000000002
$ws.Range("M73").Value = 9
$ws.Range("N73").Formula = "=SUM(K73:M73)"
$ws.Range("O73").Formula = "=N73-4000"

# --- Row 74: Chorizo Casero ---
$ws.Range("A74").Value = $dateSerial
$ws.Range("B74").Value = 13.9
$ws.Range("C74").Value = "Chorizo Casero"
$ws.Range("D74").Value = "Despensa"
$ws.Range("E74").Value = "Gasto"
$ws.Range("F74").Value = "Tarjeta Banamex"
$ws.Range("G74").Value = "Soriana"
$ws.Range("K74").Formula = "=K73-B74"
$ws.Range("L74").Value = 2527.5700000000002
$ws.Range("M74").Value = 9
$ws.Range("N74").Formula = "=SUM(K74:M74)"
$ws.Range("O74").Formula = "=N74-4000"

# --- Row 75: Pasta Dental Colgate ---
$ws.Range("A75").Value = $dateSerial
$ws.Range("B75").Value = 31.25
$ws.Range("C75").Value = "Pasta Dental Colgate"
$ws.Range("D75").Value = "Higiene"
$ws.Range("E75").Value = "Gasto"
$ws.Range("F75").Value = "Tarjeta Banamex"
$ws.Range("G75").Value = "Soriana"
$ws.Range("K75").Formula = "=K74-B75"
$ws.Range("L75").Value = 2527.5700000000002
$ws.Range("M75").Value = 9
$ws.Range("N75").Formula = "=SUM(K75:M75)"
$ws.Range("O75").Formula = "=N75-4000"

# --- Row 76: Carne de Res ---
$ws.Range("A76").Value = $dateSerial
$ws.Range("B76").Value = 43.28
$ws.Range("C76").Value = "Carne de Res"
$ws.Range("D76").Value = "Despensa"
$ws.Range("E76").Value = "Gasto"
$ws.Range("F76").Value = "Tarjeta Banamex"
$ws.Range("G76").Value = "Soriana"
$ws.Range("K76").Formula = "=K75-B76"
$ws.Range("L76").Value = 2527.5700000000002
$ws.Range("M76").Value = 9
$ws.Range("N76").Formula = "=SUM(K76:M76)"
$ws.Range("O76").Formula = "=N76-4000"

# --- Row 77: Lata de Verduras Herdez ---
$ws.Range("A77").Value = $dateSerial
$ws.Range("B77").Value = 13
$ws.Range("C77").Value = "Lata de Verduras Herdez"
$ws.Range("D77").Value = "Despensa"
$ws.Range("E77").Value = "Gasto"
$ws.Range("F77").Value = "Tarjeta Banamex"
$ws.Range("G77").Value = "Soriana"
$ws.Range("K77").Formula = "=K76-B77"
$ws.Range("L77").Value = 2527.5700000000002
$ws.Range("M77").Value = 9
$ws.Range("N77").Formula = "=SUM(K77:M77)"
$ws.Range("O77").Formula = "=N77-4000"

# --- Row 78: Frijoles con Chorizo ---
$ws.Range("A78").Value = $dateSerial
$ws.Range("B78").Value = 17
$ws.Range("C78").Value = "Frijoles con Chorizo"
$ws.Range("D78").Value = "Despensa"
$ws.Range("E78").Value = "Gasto"
$ws.Range("F78").Value = "Tarjeta Banamex"
$ws.Range("G78").Value = "Soriana"
$ws.Range("K78").Formula = "=K77-B78"
$ws.Range("L78").Value = 2527.5700000000002
$ws.Range("M78").Value = 9
$ws.Range("N78").Formula = "=SUM(K78:M78)"
$ws.Range("O78").Formula = "=N78-4000"

# --- Row 79: Huevo San Juan ---
$ws.Range("A79").Value = $dateSerial
$ws.Range("B79").Value = 25.5
$ws.Range("C79").Value = "Huevo San Juan"
$ws.Range("D79").Value = "Despensa"
$ws.Range("E79").Value = "Gasto"
$ws.Range("F79").Value = "Tarjeta Banamex"
$ws.Range("G79").Value = "Soriana"
$ws.Range("K79").Formula = "=K78-B79"
$ws.Range("L79").Value = 2527.5700000000002
$ws.Range("M79").Value = 9
$ws.Range("N79").Formula = "=SUM(K79:M79)"
$ws.Range("O79").Formula = "=N79-4000"

# --- Row 80: Leche Santa Clara ---
$ws.Range("A80").Value = $dateSerial
$ws.Range("B80").Value = 22.25
$ws.Range("C80").Value = "Leche Santa Clara"
$ws.Range("D80").Value = "Despensa"
$ws.Range("E80").Value = "Gasto"
$ws.Range("F80").Value = "Tarjeta Banamex"
$ws.Range("G80").Value = "Soriana"
$ws.Range("K80").Formula = "=K79-B80"
$ws.Range("L80").Value = 2527.5700000000002
$ws.Range("M80").Value = 9
$ws.Range("N80").Formula = "=SUM(K80:M80)"
$ws.Range("O80").Formula = "=N80-4000"

# --- Row 81: Café Soluble Nescafe ---
$ws.Range("A81").Value = $dateSerial
$ws.Range("B81").Value = 79.9
$ws.Range("C81").Value = "Café Soluble Nescafe"
$ws.Range("D81").Value = "Despensa"
$ws.Range("E81").Value = "Gasto"
$ws.Range("F81").Value = "Tarjeta Banamex"
$ws.Range("G81").Value = "Soriana"
$ws.Range("K81").Formula = "=K80-B81"
$ws.Range("L81").Value = 2527.5700000000002
$ws.Range("M81").Value = 9
$ws.Range("N81").Formula = "=SUM(K81:M81)"
$ws.Range("O81").Formula = "=N81-4000"

# --- Row 82: Sopa La Moderna ---
$ws.Range("A82").Value = $dateSerial
$ws.Range("B82").Value = 5.6
$ws.Range("C82").Value = "Sopa La Moderna"
$ws.Range("D82").Value = "Despensa"
$ws.Range("E82").Value = "Gasto"
$ws.Range("F82").Value = "Tarjeta Banamex"
$ws.Range("G82").Value = "Soriana"
$ws.Range("K82").Formula = "=K81-B82"
$ws.Range("L82").Value = 2527.5700000000002
$ws.Range("M82").Value = 9
$ws.Range("N82").Formula = "=SUM(K82:M82)"
$ws.Range("O82").Formula = "=N82-4000"

# --- Row 83: Sopa La Moderna ---
$ws.Range("A83").Value = $dateSerial
$ws.Range("B83").Value = 5.6
$ws.Range("C83").Value = "Sopa La Moderna"
$ws.Range("D83").Value = "Despensa"
$ws.Range("E83").Value = "Gasto"
$ws.Range("F83").Value = "Tarjeta Banamex"
$ws.Range("G83").Value = "Soriana"
$ws.Range("K83").Formula = "=K82-B83"
$ws.Range("L83").Value = 2527.5700000000002
$ws.Range("M83").Value = 9
$ws.Range("N83").Formula = "=SUM(K83:M83)"
$ws.Range("O83").Formula = "=N83-4000"

# --- Row 84: Tortilla de Maiz ---
$ws.Range("A84").Value = $dateSerial
$ws.Range("B84").Value = 24.75
$ws.Range("C84").Value = "Tortilla de Maiz"
$ws.Range("D84").Value = "Despensa"
$ws.Range("E84").Value = "Gasto"
$ws.Range("F84").Value = "Tarjeta Banamex"
$ws.Range("G84").Value = "Soriana"
$ws.Range("K84").Formula = "=K83-B84"
$ws.Range("L84").Value = 2527.5700000000002
$ws.Range("M84").Value = 9
$ws.Range("N84").Formula = "=SUM(K84:M84)"
$ws.Range("O84").Formula = "=N84-4000"

# --- Row 85: Propina ---
$ws.Range("A85").Value = $dateSerial
$ws.Range("B85").Value = 7
$ws.Range("C85").Value = "Propina"
$ws.Range("D85").Value = "Misc"
$ws.Range("E85").Value = "Gasto"
$ws.Range("F85").Value = "Efectivo"
$ws.Range("G85").Value = "Soriana"
$ws.Range("K85").Value = 7601.91
$ws.Range("L85").Value = 2527.5700000000002
$ws.Range("M85").Formula = "=M84-B85"
$ws.Range("N85").Formula = "=SUM(K85:M85)"
$ws.Range("O85").Formula = "=N85-4000"

# --- Row 86: Plan AT&T ---
$ws.Range("A86").Value = $dateSerial
$ws.Range("B86").Value = 240.07
$ws.Range("C86").Value = "Plan AT&T"
$ws.Range("D86").Value = "Servicios"
$ws.Range("E86").Value = "Gasto"
$ws.Range("F86").Value = "Tarjeta Banamex"
$ws.Range("G86").Value = "AT&T"
$ws.Range("K86").Formula = "=K85-B86"
$ws.Range("L86").Value = 2527.5700000000002
$ws.Range("M86").Value = 2
$ws.Range("N86").Formula = "=SUM(K86:M86)"
$ws.Range("O86").Formula = "=N86-4000"

# Match the date-number formatting used by the rest of column A (style copied from A71).
$ws.Range("A71").Copy()
$ws.Range("A72:A86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 83 duplicates row 82 formatting (B:G), matching the original author's copy/paste.
$ws.Range("B82:G82").Copy()
$ws.Range("B83:G83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# N/O total columns on "boundary" rows pick up the same style as the rest of the totals column.
$ws.Range("N61:O61").Copy()
$ws.Range("N72:O72").PasteSpecial(-4122)
$ws.Range("N61:O61").Copy()
$ws.Range("N73:O73").PasteSpecial(-4122)
$ws.Range("N61:O61").Copy()
$ws.Range("N74:O74").PasteSpecial(-4122)
$ws.Range("N61:O61").Copy()
$ws.Range("N75:O75").PasteSpecial(-4122)
$ws.Range("N61:O61").Copy()
$ws.Range("N85:O85").PasteSpecial(-4122)
$ws.Range("N61:O61").Copy()
$ws.Range("N86:O86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the author left off.
$ws.Range("R70").Select()

